# Auto-generated edit script: apply Seraph_Profits.xlsx diff to before.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 188.85715
$ws.Range("I11").Value = 188.85715
$ws.Range("K11").Value = 188.85715
$ws.Range("M11").Value = -48.85714999999999

# Row 70
$ws.Range("H70").Value = 80323.09
$ws.Range("J70").Value = 106695.25
$ws.Range("L70").Value = 320085.75
$ws.Range("N70").Value = -320625.75

# Row 73
$ws.Range("H73").Value = 80323.09
$ws.Range("J73").Value = 106695.25
$ws.Range("L73").Value = 320085.75
$ws.Range("N73").Value = -321957.75

# Row 80
$ws.Range("H80").Value = 21943.572
$ws.Range("I80").Value = 536.6667
$ws.Range("J80").Value = 37998.75
$ws.Range("K80").Value = 1610.0001
$ws.Range("L80").Value = 113996.25
$ws.Range("M80").Value = -612.0001
$ws.Range("N80").Value = -115992.25

# Row 83
$ws.Range("H83").Value = 21943.572
$ws.Range("I83").Value = 536.6667
$ws.Range("J83").Value = 37998.75
$ws.Range("K83").Value = 4830.0003
$ws.Range("L83").Value = 341988.75
$ws.Range("M83").Value = 161.9997000000003
$ws.Range("N83").Value = -351972.75

# Row 113
$ws.Range("H113").Value = 3634.4443
$ws.Range("I113").Value = 3644.2856
$ws.Range("K113").Value = 3644.2856
$ws.Range("M113").Value = -390.2856000000002

# Row 132
$ws.Range("H132").Value = 2716.3333
$ws.Range("I132").Value = 2710
$ws.Range("K132").Value = 8130
$ws.Range("M132").Value = -5600

# Row 138
$ws.Range("H138").Value = 7715.591
$ws.Range("I138").Value = 4538.6
$ws.Range("J138").Value = 8650
$ws.Range("K138").Value = 13615.8
$ws.Range("L138").Value = 25950
$ws.Range("M138").Value = -8475.800000000001
$ws.Range("N138").Value = -36230

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1921.2
$ws.Range("I45").Value = 1901.5
$ws.Range("K45").Value = 1901.5
$ws.Range("M45").Value = -1524.5

# Row 61
$ws.Range("H61").Value = 1554.3334
$ws.Range("I61").Value = 1458.8462
$ws.Range("K61").Value = 1458.8462
$ws.Range("M61").Value = -1246.8462

# Row 97
$ws.Range("H97").Value = 746.6
$ws.Range("I97").Value = 746.6
$ws.Range("K97").Value = 746.6
$ws.Range("M97").Value = -250.6

# Row 102
$ws.Range("H102").Value = 1658
$ws.Range("I102").Value = 1658
$ws.Range("K102").Value = 1658
$ws.Range("M102").Value = -36

# Row 132
$ws.Range("H132").Value = 1850.159
$ws.Range("I132").Value = 1640.375
$ws.Range("K132").Value = 4921.125
$ws.Range("M132").Value = -2391.125

# Row 136
$ws.Range("H136").Value = 1554.3334
$ws.Range("I136").Value = 1458.8462
$ws.Range("K136").Value = 4376.5386
$ws.Range("M136").Value = -1826.5386

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1966.6666
$ws.Range("I20").Value = 1950
$ws.Range("K20").Value = 1950
$ws.Range("M20").Value = -1703

# Row 125
$ws.Range("H125").Value = 30000
$ws.Range("J125").Value = 30000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -39840

# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 8
$ws.Range("H8").Value = 5274.5
$ws.Range("I8").Value = 5274.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 5274.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -5134.5
$ws.Range("N8").ClearContents()

# Row 31
$ws.Range("H31").Value = 4650.2104
$ws.Range("I31").Value = 2435.375
$ws.Range("J31").Value = 6261
$ws.Range("K31").Value = 2435.375
$ws.Range("L31").Value = 6261
$ws.Range("M31").Value = -2140.375
$ws.Range("N31").Value = -6851

# Row 34
$ws.Range("H34").Value = 4650.2104
$ws.Range("I34").Value = 2435.375
$ws.Range("J34").Value = 6261
$ws.Range("K34").Value = 2435.375
$ws.Range("L34").Value = 6261
$ws.Range("M34").Value = -2233.375
$ws.Range("N34").Value = -6665

# Row 58
$ws.Range("H58").Value = 5459.5557
$ws.Range("J58").Value = 6934.2
$ws.Range("L58").Value = 6934.2
$ws.Range("N58").Value = -7340.2

# Row 62
$ws.Range("H62").Value = 4037.6924
$ws.Range("I62").Value = 3799.9
$ws.Range("J62").Value = 4830.3335
$ws.Range("K62").Value = 3799.9
$ws.Range("L62").Value = 4830.3335
$ws.Range("M62").Value = -3175.9
$ws.Range("N62").Value = -6078.3335

# Row 65
$ws.Range("H65").Value = 4037.6924
$ws.Range("I65").Value = 3799.9
$ws.Range("J65").Value = 4830.3335
$ws.Range("K65").Value = 18999.5
$ws.Range("L65").Value = 24151.6675
$ws.Range("M65").Value = -15879.5
$ws.Range("N65").Value = -30391.6675

# Row 80
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

# Row 83
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

# Row 122
$ws.Range("H122").Value = 6869.3076
$ws.Range("I122").Value = 6872.3335
$ws.Range("J122").Value = 6862.5
$ws.Range("K122").Value = 20617.0005
$ws.Range("L122").Value = 20587.5
$ws.Range("M122").Value = -18167.0005
$ws.Range("N122").Value = -25487.5

# Row 136
$ws.Range("H136").Value = 5459.5557
$ws.Range("J136").Value = 6934.2
$ws.Range("L136").Value = 20802.6
$ws.Range("N136").Value = -25902.6

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1615.091
$ws.Range("J68").Value = 2063.6
$ws.Range("L68").Value = 6190.799999999999
$ws.Range("N68").Value = -7812.799999999999

# Row 71
$ws.Range("H71").Value = 1615.091
$ws.Range("J71").Value = 2063.6
$ws.Range("L71").Value = 18572.4
$ws.Range("N71").Value = -26684.4

# Row 106
$ws.Range("H106").Value = 20987
$ws.Range("J106").Value = 20987
$ws.Range("L106").Value = 62961
$ws.Range("N106").Value = -64853

# Row 131
$ws.Range("H131").Value = 3156.25
$ws.Range("J131").Value = 3890
$ws.Range("L131").Value = 11670
$ws.Range("N131").Value = -21750

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 17500
$ws.Range("J15").Value = 17500
$ws.Range("L15").Value = 17500
$ws.Range("N15").Value = -18076

# Row 81
$ws.Range("H81").Value = 17500
$ws.Range("J81").Value = 17500
$ws.Range("L81").Value = 17500
$ws.Range("N81").Value = -19496

# Row 84
$ws.Range("H84").Value = 17500
$ws.Range("J84").Value = 17500
$ws.Range("L84").Value = 52500
$ws.Range("N84").Value = -62484

# Row 94
$ws.Range("H94").Value = 13333.333
$ws.Range("J94").Value = 13333.333
$ws.Range("L94").Value = 13333.333
$ws.Range("N94").Value = -14685.333

# Row 102
$ws.Range("H102").Value = 1157.1708
$ws.Range("I102").Value = 614.85187
$ws.Range("J102").Value = 2203.0715
$ws.Range("K102").Value = 614.85187
$ws.Range("L102").Value = 2203.0715
$ws.Range("M102").Value = 1007.14813
$ws.Range("N102").Value = -5447.0715

# Row 132
$ws.Range("H132").Value = 2723.1052
$ws.Range("I132").Value = 2123.4
$ws.Range("J132").Value = 4972
$ws.Range("K132").Value = 6370.200000000001
$ws.Range("L132").Value = 14916
$ws.Range("M132").Value = -3840.200000000001
$ws.Range("N132").Value = -19976

# Row 136
$ws.Range("H136").Value = 29886.111
$ws.Range("J136").Value = 29886.111
$ws.Range("L136").Value = 89658.333
$ws.Range("N136").Value = -94758.333

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2089.4546
$ws.Range("I40").Value = 1748.5
$ws.Range("J40").Value = 2998.6667
$ws.Range("K40").Value = 1748.5
$ws.Range("L40").Value = 2998.6667
$ws.Range("M40").Value = -1612.5
$ws.Range("N40").Value = -3270.6667

# Row 80
$ws.Range("H80").Value = 33000
$ws.Range("J80").Value = 33000
$ws.Range("L80").Value = 33000
$ws.Range("N80").Value = -35246

# Row 83
$ws.Range("H83").Value = 33000
$ws.Range("J83").Value = 33000
$ws.Range("L83").Value = 99000
$ws.Range("N83").Value = -110232

# Row 100
$ws.Range("H100").Value = 2634.3333
$ws.Range("I100").Value = 2451.5
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2451.5
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1910.5
$ws.Range("N100").Value = -4082

# Row 104
$ws.Range("H104").Value = 37499.75
$ws.Range("J104").Value = 37499.75
$ws.Range("L104").Value = 37499.75
$ws.Range("N104").Value = -44487.75

# Row 122
$ws.Range("H122").Value = 3510
$ws.Range("I122").Value = 3483.5
$ws.Range("K122").Value = 10450.5
$ws.Range("M122").Value = -8000.5

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 7376.357
$ws.Range("I62").Value = 5156.3335
$ws.Range("K62").Value = 5156.3335
$ws.Range("M62").Value = -4532.3335

# Row 65
$ws.Range("H65").Value = 7376.357
$ws.Range("I65").Value = 5156.3335
$ws.Range("K65").Value = 25781.6675
$ws.Range("M65").Value = -22661.6675

# Row 126
$ws.Range("H126").Value = 79530.08
$ws.Range("I126").Value = 112488.445
$ws.Range("J126").Value = 5373.75
$ws.Range("K126").Value = 337465.335
$ws.Range("L126").Value = 16121.25
$ws.Range("M126").Value = -334995.335
$ws.Range("N126").Value = -21061.25

# Row 132
$ws.Range("H132").Value = 2158.25
$ws.Range("I132").Value = 1790
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 5370
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -2840
$ws.Range("N132").Value = -17058.5
